# CodeSystem-BRImunobiologico.xlsx
# - Rename concept "100" (row 101) from VVBN / Vacina Varíola Bavarian Nordic
#   to VVS / Vacina Varíola Símia (Atenuada)
# - Add two new concepts at the end of the "Concepts" table:
#     113 CHIKUNGUNYA        / Vacina Chikungunya (recombinante e atenuada)
#     114 COVID-19 SINOPHARM / Vacina Covid-19-inativada, Sinopharm

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# --- 1) Update existing row 101 (code "100") ---------------------------
$ws.Cells.Item(101, 3).Value = "VVS"
$ws.Cells.Item(101, 4).Value = "Vacina Varíola Símia (Atenuada)"

# --- 2) Append two new data rows, matching the look of the last row ----
$ws.Range("A113:D113").Copy()
$ws.Range("A114:D115").PasteSpecial(-4122)   # xlPasteFormats

# Helper pattern: write a numeric-looking code ("1", "113", "114") as
# genuine text (not a number) without leaving stray number-format/style
# entries behind - round-trip it through a throw-away text formula cell
# and paste only its (string) value.

$ws.Cells.Item(500, 1).Formula = "=""1"""
$ws.Cells.Item(500, 1).Copy()
$ws.Cells.Item(114, 1).PasteSpecial(-4163)   # xlPasteValues

$ws.Cells.Item(500, 1).Formula = "=""113"""
$ws.Cells.Item(500, 1).Copy()
$ws.Cells.Item(114, 2).PasteSpecial(-4163)

$ws.Cells.Item(114, 3).Value = "CHIKUNGUNYA"
$ws.Cells.Item(114, 4).Value = "Vacina Chikungunya (recombinante e atenuada)"

$ws.Cells.Item(500, 1).Formula = "=""1"""
$ws.Cells.Item(500, 1).Copy()
$ws.Cells.Item(115, 1).PasteSpecial(-4163)

$ws.Cells.Item(500, 1).Formula = "=""114"""
$ws.Cells.Item(500, 1).Copy()
$ws.Cells.Item(115, 2).PasteSpecial(-4163)

$ws.Cells.Item(115, 3).Value = "COVID-19 SINOPHARM"
$ws.Cells.Item(115, 4).Value = "Vacina Covid-19-inativada, Sinopharm"

# Clean up the scratch cell used for the text round-trips above
$ws.Cells.Item(500, 1).Delete()
